# Add Animation and HeathBar for one character.
# - Adds two new worksheets ("Battle", "Effect") with Wukong skill / status-effect
#   localization rows, appends 6 new rows to the existing "STR" sheet, and leaves
#   the STR sheet as the active / selected tab.

$wb = $excel.ActiveWorkbook

$strSheet = $wb.Worksheets.Item("STR")

# --- Create the two new sheets, appended after STR (last sheet) -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battle = $wb.Worksheets.Add($null, $lastSheet)
$battle.Name = "Battle"
$effect = $wb.Worksheets.Add($null, $battle)
$effect.Name = "Effect"

# --- Header formatting (copy the ID/ENGLISH/VIETNAMESE header style from STR) ----
# PasteSpecial(formats) must run before the header values are written, because it
# blanks out whatever is currently in the destination cells.
$strSheet.Range("A1:C1").Copy() | Out-Null
$battle.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$effect.Range("A1:C1").PasteSpecial(-4122) | Out-Null

$battle.Range("A1").Value = "ID"
$battle.Range("B1").Value = "ENGLISH"
$battle.Range("C1").Value = "VIETNAMESE"

$effect.Range("A1").Value = "ID"
$effect.Range("B1").Value = "ENGLISH"
$effect.Range("C1").Value = "VIETNAMESE"

# --- Battle sheet: Wukong skill name/description rows ----------------------------
$battle.Range("A2").Value = "Wukong_FirstSkil_Name"
$battle.Range("C2").Value = "Đòn Tàn Khốc"
$battle.Range("A3").Value = "Wukong_FirstSkil_Des"
$battle.Range("C3").Value = "Gây ST bằng 80% Tấn Công cho 1 kẻ địch"
$battle.Range("A4").Value = "Wukong_SecondSkill_Name"
$battle.Range("C4").Value = "Phân Thân"
$battle.Range("A5").Value = "Wukong_SecondSkill_Des"
$battle.Range("C5").Value = "Gây 60% ST Tấn Công cho toàn bộ kẻ định đồng thời hồi 30% ST gây ra."
$battle.Range("A6").Value = "Wukong_ThirdSkill_Name"
$battle.Range("C6").Value = "Thiết Bảng Ngàn Cân"
$battle.Range("A7").Value = "Wukong_ThirdSkill_Des"
$battle.Range("C7").Value = "Gây 180% ST Tấn Công cho kẻ định, kèm thêm Giảm Phòng Thù cho mục tiêu, duy trì 2 hiệp"

# --- Effect sheet: status-effect name/description rows ---------------------------
$effect.Range("A2").Value = "DEBUFF_DEF_NAME"
$effect.Range("C2").Value = "Giảm Phòng Thủ"
$effect.Range("A3").Value = "DEBUFF_DEF_DES"
$effect.Range("C3").Value = "Giảm Phòng Thủ {0}%."
$effect.Range("A4").Value = "RESIST_DEBUFF_NAME"
$effect.Range("C4").Value = "Miễn Dịch"
$effect.Range("A5").Value = "RESIST_DEBUFF_DES"
$effect.Range("C5").Value = "Không bị ảnh hưởng bởi tất cả các hiệu ứng suy yếu."
$effect.Range("A6").Value = "BUFF_DEF_NAME"
$effect.Range("C6").Value = "Tăng Phòng Thủ"
$effect.Range("A7").Value = "BUFF_DEF_DES"
$effect.Range("C7").Value = "Phòng Thủ Tăng {0}%."
$effect.Range("A8").Value = "DEBUFF_RECEIVE_NAME"
$effect.Range("C8").Value = "Không thể cường hóa"
$effect.Range("A9").Value = "DEBUFF_RECEIVE_DES"
$effect.Range("C9").Value = "Không thể nhận được bất kì hiệu ứng cường hóa."
$effect.Range("A10").Value = "BUFF_ATTACK_NAME"
$effect.Range("C10").Value = "Tăng Tấn Công"
$effect.Range("A11").Value = "BUFF_ATTACK_DES"
$effect.Range("C11").Value = "Tấn Công tăng {0}%."
$effect.Range("A12").Value = "DEBUFF_ATTACK_NAME"
$effect.Range("C12").Value = "Giảm Tấn Công"
$effect.Range("A13").Value = "DEBUFF_ATTACK_DES"
$effect.Range("C13").Value = "Tấn Công giảm {0}%."
$effect.Range("A14").Value = "DEBUFF_HEAL_NAME"
$effect.Range("C14").Value = "Cấm Điều Trị"
$effect.Range("A15").Value = "DEBUFF_HEAL_DES"
$effect.Range("C15").Value = "Không thể nhận hiệu ứng hồi máu."
$effect.Range("A16").Value = "DEBUFF_STUN_NAME"
$effect.Range("C16").Value = "Choáng"
$effect.Range("A17").Value = "DEBUFF_STUN_DES"
$effect.Range("C17").Value = "Không thể hành động trong thời gian duy trì hiệp."

# --- STR sheet: append 6 new rows (167-172) for generic skill-level strings -------
$strSheet.Range("A167").Value = "STR_SLILL_LEVEL_NAME"
$strSheet.Range("C167").Value = "Nhân vật Lv.{0}"
$strSheet.Range("A168").Value = "STR_ST_SKIL_DES"
$strSheet.Range("C168").Value = "ST tăng tới {0} Tấn Công."
$strSheet.Range("A169").Value = "STR_CRIT_COOLDOWN_ULTIMATE"
$strSheet.Range("C169").Value = "Nếu kĩ năng này Bạo Kích, hiệp hồi chiêu của {0} -1."
$strSheet.Range("A170").Value = "STR_DMG_CRIT"
$strSheet.Range("C170").Value = "Nếu kĩ năng này Bạo Kích, ST tăng {0}."
$strSheet.Range("A171").Value = "STR_SKILL_COOLDOWN"
$strSheet.Range("C171").Value = "Hiệp hồi chiêu -1."
$strSheet.Range("A172").Value = "STR_RESIST_DEFBUFF"
$strSheet.Range("C172").Value = "Banr thân được Miễn Dịch, duy trì 2 hiệp."

# --- Column widths (approximate auto-fit sizing for the new sheets) --------------
$battle.Columns.Item(1).ColumnWidth = 33.5
$battle.Columns.Item(2).ColumnWidth = 30.833333333333336
$battle.Columns.Item(3).ColumnWidth = 35.83333333333333

$effect.Columns.Item(1).ColumnWidth = 17.833333333333336
$effect.Columns.Item(3).ColumnWidth = 17.666666666666668

# --- Selections on each touched/new sheet -----------------------------------------
$battle.Range("A1:C1").Select() | Out-Null
$effect.Range("A17").Select() | Out-Null

# STR keeps its zoom, scrolls further down, and selects the last new cell - it is
# also the sheet left active/selected when the workbook is saved.
$strSheet.Range("C170").Select() | Out-Null
$strSheet.Activate() | Out-Null
